# Apply updated crypto price/volume data per Mon Apr 22 18:20:42 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.544.36"
$ws.Range("E2").Value = "  +2.67%  "

$ws.Range("D3").Value = "3.200.69"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'597.81"
$ws.Range("E5").Value = "  +3.69%  "

$ws.Range("D6").Value = "'155.04"
$ws.Range("E6").Value = "  +4.16%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.199.72"
$ws.Range("E8").Value = "  +1.73%  "

$ws.Range("E9").Value = "  +4.08%  "

$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("D11").Value = "'6.02"
$ws.Range("E11").Value = "  -0.88%  "

$ws.Range("D12").Value = "'0.517"
$ws.Range("E12").Value = "  +4.33%  "

$ws.Range("E13").Value = "  +4.08%  "

$ws.Range("D14").Value = "'39.13"
$ws.Range("E14").Value = "  +6.04%  "

$ws.Range("D15").Value = "3.726.58"
$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("D16").Value = "66.533.50"
$ws.Range("E16").Value = "  +2.46%  "

$ws.Range("E17").Value = "  +5.58%  "

$ws.Range("D18").Value = "3.202.35"
$ws.Range("E18").Value = "  +1.84%  "

$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").Value = "'513.54"
$ws.Range("E20").Value = "  +2.18%  "

$ws.Range("E21").Value = "  +4.42%  "

$ws.Range("D22").Value = "'0.742"
$ws.Range("E22").Value = "  +4.64%  "

$ws.Range("E23").Value = "  +5.65%  "

$ws.Range("E24").Value = "  -0.48%  "

$ws.Range("D25").Value = "'85.79"
$ws.Range("E25").Value = "  +2.51%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.28%  "

$ws.Range("D27").Value = "'9.32"
$ws.Range("E27").Value = "  +5.80%  "

$ws.Range("D28").Value = "'3.02"
$ws.Range("E28").Value = "  +4.62%  "

$ws.Range("E29").Value = "  +7.12%  "

$ws.Range("D30").Value = "'7.14"
$ws.Range("E30").Value = "  +16.12%  "

$ws.Range("D31").Value = "'2.94"
$ws.Range("E31").Value = "  +5.27%  "

$ws.Range("D32").Value = "'28.30"
$ws.Range("E32").Value = "  +3.25%  "

$ws.Range("D33").Value = "'1.24"
$ws.Range("E33").Value = "  +4.16%  "

$ws.Range("E34").Value = "  +0.20%  "

$ws.Range("D35").Value = "'6.55"
$ws.Range("E35").Value = "  +1.90%  "

$ws.Range("D36").Value = "'522.89"
$ws.Range("E36").Value = "  +9.35%  "

$ws.Range("D37").Value = "'54.87"
$ws.Range("E37").Value = "  +0.55%  "

$ws.Range("D38").Value = "'0.0899"
$ws.Range("E38").Value = "  +1.21%  "

$ws.Range("D39").Value = "'0.0424"
$ws.Range("E39").Value = "  +2.42%  "

$ws.Range("D40").Value = "'8.90"
$ws.Range("E40").Value = "  +3.61%  "

$ws.Range("E41").Value = "  +6.72%  "

$ws.Range("D42").Value = "0.0₃0683"
$ws.Range("E42").Value = "  +18.73%  "

$ws.Range("E43").Value = "  -2.43%  "

$ws.Range("D44").Value = "'0.303"
$ws.Range("E44").Value = "  +7.84%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.46"
$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.931.96"
$ws.Range("E46").Value = "  -2.32%  "

$ws.Range("D47").Value = "'28.68"
$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("E48").Value = "  +3.28%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").Value = "'2.33"
$ws.Range("E50").Value = "  +4.91%  "

$ws.Range("D51").Value = "'2.63"
$ws.Range("E51").Value = "  +10.88%  "

